$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.505.59'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.922.10'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.60'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4846'
$ws.Range('E7').Value = '  +3.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4081'
$ws.Range('E8').Value = '  +1.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08180'
$ws.Range('E9').Value = '  +2.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.018'
$ws.Range('E10').Value = '  +2.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.76'
$ws.Range('E11').Value = '  +5.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.928.34'
$ws.Range('E12').Value = '  +1.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.049'
$ws.Range('E13').Value = '  +3.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.241'
$ws.Range('E14').Value = '  +3.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.25'
$ws.Range('E15').Value = '  +2.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06765'
$ws.Range('E16').Value = '  +2.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.005'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('E18').Value = '  +1.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.78'
$ws.Range('E19').Value = '  +2.13%  '
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '29.526.60'
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.642'
$ws.Range('E22').Value = '  +2.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.75'
$ws.Range('E23').Value = '  +2.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.178'
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.156.22'
$ws.Range('E25').Value = '  +1.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.638'
$ws.Range('E26').Value = '  +11.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.68'
$ws.Range('E27').Value = '  +1.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.12'
$ws.Range('E28').Value = '  +2.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.132'
$ws.Range('E29').Value = '  +2.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.79'
$ws.Range('E30').Value = '  +3.49%  '
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09554'
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.532'
$ws.Range('E33').Value = '  +3.81%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.399'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.554'
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02286'
$ws.Range('E36').Value = '  +2.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06133'
$ws.Range('E37').Value = '  +1.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.188'
$ws.Range('E38').Value = '  +1.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '10.92'
$ws.Range('E39').Value = '  +8.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5990'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.032'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1861'
$ws.Range('E42').Value = '  +2.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.430'
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.279'
$ws.Range('E44').Value = '  +2.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.07634'
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.42'
$ws.Range('E46').Value = '  +2.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5595'
$ws.Range('E47').Value = '  +2.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.960'
$ws.Range('E48').Value = '  +3.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '116.70'
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.94'
$ws.Range('E50').Value = '  +3.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.423'
$ws.Range('E51').Value = '  +3.92%  '
